# Applies the "Completed with exception handling and annotations" edit:
#  - Settings!C6  : description for OutputSheet
#  - Constants!C18: description for BrowserPath
#  - Constants!A19:C19 : new SystemException row
#  - Constants!A20:C20 : new TimeOut (delay) row
#  - Constants!A21:C21 : new BusinessException row

$wb = $excel.ActiveWorkbook

# --- Settings sheet -------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("C6").Value = "Sheet where the output data is stored"

# --- Constants sheet -------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

# Description for the existing BrowserPath row
$constants.Range("C18").Value = "Botdna notaries url"

# New row 19 - SystemException / Page not loading message
$constants.Range("A19").Value = "SystemException"
$constants.Range("B19").Value = "Page not  loading"
$constants.Range("C19").Value = "Message to be displayedwhen the system exception occurs"

# New row 20 - TimeOut delay for check app state activity
$constants.Range("A20").Value = "TimeOut"
$constants.Range("B20").Value = 5
$constants.Range("C20").Value = "Delay for the check app state activity"

# New row 21 - BusinessException / incorrect details message
$constants.Range("A21").Value = "BusinessException"
$constants.Range("B21").Value = "Enetered details are incorrect"
$constants.Range("C21").Value = "Message to be displayedwhen the Business exception occurs"

# --- Row heights reflow slightly once the wrapped description text ---
# --- is re-measured by the editing session (matches the saved file) --
$settings.Rows.Item(3).RowHeight = 14.5
$settings.Rows.Item(5).RowHeight = 29
$constants.Rows.Item(2).RowHeight = 29
$constants.Rows.Item(3).RowHeight = 43.5
$constants.Rows.Item(17).RowHeight = 29

# --- Restore the final selections shown in the saved workbook --------
# (Settings stays the active/front tab, so select there last.)
$constants.Activate()
$constants.Range("C22").Select()
$settings.Activate()
$settings.Range("C7").Select()
